$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New id/label pairs for rows 2..40 (B and C columns), with D (is_prefered) cleared.
$rows = @(
    @{B='#tvveeden-soldaet'; C='Tvveeden Soldaet'},
    @{B='#caesar'; C='Caesar'},
    @{B='#charon'; C='Charon'},
    @{B='#eros'; C='Eros'},
    @{B='#m.-agrippa'; C='M. Agrippa'},
    @{B='#c.-augustus'; C='C. Augustus'},
    @{B='#dercetaeus'; C='Dercetaeus'},
    @{B="#anthonius-in't-harnas"; C="Anthonius in't harnas"},
    @{B='#eerste-slaef'; C='Eerste slaef'},
    @{B='#caesar-augustus'; C='Caesar Augustus'},
    @{B='#charmion'; C='Charmion'},
    @{B='#proculeius'; C='Proculeius'},
    @{B='#thyreus'; C='Thyreus'},
    @{B='#landt-man'; C='Landt-man'},
    @{B='#derde-slaef'; C='Derde slaef'},
    @{B='#plancus'; C='Plancus'},
    @{B='#lucillus'; C='Lucillus'},
    @{B='#tvveede-slaef'; C='Tvveede slaef'},
    @{B='#gallus'; C='Gallus'},
    @{B='#caluisius'; C='Caluisius'},
    @{B='#agrippa'; C='Agrippa'},
    @{B='#anthonius'; C='Anthonius'},
    @{B='#niger'; C='Niger'},
    @{B='#titus'; C='Titus'},
    @{B='#hopman'; C='Hopman'},
    @{B='#euphronius'; C='Euphronius'},
    @{B='#ostrobas'; C='Ostrobas'},
    @{B='#aristocrates'; C='Aristocrates'},
    @{B='#cleopatra'; C='Cleopatra'},
    @{B='#eersten-soldaet'; C='Eersten Soldaet'},
    @{B='#soldaat'; C='Soldaat'},
    @{B='#daemon'; C='Daemon'},
    @{B='#soldaet'; C='Soldaet'},
    @{B='#diomedes'; C='Diomedes'},
    @{B='#seleucus'; C='Seleucus'},
    @{B='#canidius'; C='Canidius'},
    @{B='#iras'; C='Iras'},
    @{B='#anthoniado'; C='Anthoniado'},
    @{B='#vierde-slaef'; C='Vierde slaef'}
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value = $rows[$i].B
    $ws.Cells.Item($r, 3).Value = $rows[$i].C
    $ws.Cells.Item($r, 4).Value = ""
}
